$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "0.999", "1.40" etc. must not be auto-converted to numbers).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D12", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D43", "D45", "D46")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '61.707.14'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '3.398.30'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '577.60'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').Value = '143.30'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').Value = '3.978.67'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '28.00'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '3.395.09'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '61.709.99'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '6.15'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '13.68'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '9.15'
$ws.Range('E20').Value = '  +2.01%  '
$ws.Range('D21').Value = '389.40'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').Value = '74.67'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '0.550'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').Value = '0.0000114'
$ws.Range('E25').Value = '  -2.89%  '
$ws.Range('D26').Value = '0.182'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '7.42'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '8.01'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.40'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = '23.41'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('D34').Value = '6.95'
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').Value = '168.33'
$ws.Range('E35').Value = '  +1.26%  '
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('D37').Value = '3.430.44'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').Value = '0.0763'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').Value = '27.18'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '0.783'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('E42').Value = '  +1.27%  '
$ws.Range('D43').Value = '1.68'
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('D45').Value = '2.479.50'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('D46').Value = '22.78'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('E51').Value = '  -1.20%  '
